# Add team record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - AD, AE, AF with same bold/centered style as existing headers (copy style from AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-47: Wins=76, Losses=86, Ties=0 for every row
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD
    $ws.Cells.Item($r, 31).Value = 86   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
